$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.967.58"
Set-TextValue $ws.Range("E2") "  +0.22%  "

Set-TextValue $ws.Range("D3") "1.638.86"
Set-TextValue $ws.Range("E3") "  -0.19%  "

Set-TextValue $ws.Range("E4") "  +0.06%  "

Set-TextValue $ws.Range("D5") "212.41"
Set-TextValue $ws.Range("E5") "  -0.13%  "

Set-TextValue $ws.Range("E6") "  -0.08%  "

Set-TextValue $ws.Range("E7") "  -0.03%  "

Set-TextValue $ws.Range("D8") "23.38"
Set-TextValue $ws.Range("E8") "  -0.30%  "

Set-TextValue $ws.Range("D10") "0.0613"
Set-TextValue $ws.Range("E10") "  -0.01%  "

Set-TextValue $ws.Range("E11") "  +1.63%  "

Set-TextValue $ws.Range("E12") "  -0.10%  "

Set-TextValue $ws.Range("D13") "1.642.66"
Set-TextValue $ws.Range("E13") "  -0.28%  "

Set-TextValue $ws.Range("E14") "  +0.13%  "

Set-TextValue $ws.Range("D15") "0.569"
Set-TextValue $ws.Range("E15") "  +1.04%  "

Set-TextValue $ws.Range("D16") "65.46"

Set-TextValue $ws.Range("D17") "27.966.13"
Set-TextValue $ws.Range("E17") "  +0.32%  "

Set-TextValue $ws.Range("D18") "232.88"
Set-TextValue $ws.Range("E18") "  +0.53%  "

Set-TextValue $ws.Range("D19") "0.0₃0721"
Set-TextValue $ws.Range("E19") "  -0.34%  "

Set-TextValue $ws.Range("D20") "7.54"
Set-TextValue $ws.Range("E20") "  -1.85%  "

Set-TextValue $ws.Range("E21") "  +0.01%  "

Set-TextValue $ws.Range("D22") "10.42"
Set-TextValue $ws.Range("E22") "  -3.25%  "

Set-TextValue $ws.Range("E23") "  -0.36%  "

Set-TextValue $ws.Range("E24") "  -3.85%  "

Set-TextValue $ws.Range("D25") "153.14"

Set-TextValue $ws.Range("D26") "6.94"

Set-TextValue $ws.Range("D27") "15.65"
Set-TextValue $ws.Range("E27") "  -0.30%  "

Set-TextValue $ws.Range("E28") "  -0.76%  "

Set-TextValue $ws.Range("E29") "  +0.04%  "

Set-TextValue $ws.Range("E30") "  +0.58%  "

Set-TextValue $ws.Range("D31") "0.0485"
Set-TextValue $ws.Range("E31") "  +0.33%  "

Set-TextValue $ws.Range("D32") "3.38"
Set-TextValue $ws.Range("E32") "  +2.32%  "

Set-TextValue $ws.Range("B33") "Maker"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D33") "1.404.38"
Set-TextValue $ws.Range("E33") "  -3.42%  "

Set-TextValue $ws.Range("B34") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "3.08"
Set-TextValue $ws.Range("E34") "  -0.53%  "

Set-TextValue $ws.Range("E35") "  +1.57%  "

Set-TextValue $ws.Range("E36") "  +1.25%  "

Set-TextValue $ws.Range("E37") "  +0.29%  "

Set-TextValue $ws.Range("E38") "  -0.43%  "

Set-TextValue $ws.Range("E39") "  +1.07%  "

Set-TextValue $ws.Range("E40") "  -1.27%  "

Set-TextValue $ws.Range("E41") "  +0.70%  "

Set-TextValue $ws.Range("E42") "  -0.06%  "

Set-TextValue $ws.Range("D43") "67.13"
Set-TextValue $ws.Range("E43") "  -3.05%  "

Set-TextValue $ws.Range("E44") "  +3.08%  "

Set-TextValue $ws.Range("D45") "1.82"
Set-TextValue $ws.Range("E45") "  +2.43%  "

Set-TextValue $ws.Range("E46") "  -0.54%  "

Set-TextValue $ws.Range("D47") "1.780.75"
Set-TextValue $ws.Range("E47") "  -0.19%  "

Set-TextValue $ws.Range("D48") "88.13"
Set-TextValue $ws.Range("E48") "  -0.32%  "

Set-TextValue $ws.Range("D49") "0.100"
Set-TextValue $ws.Range("E49") "  -0.32%  "

Set-TextValue $ws.Range("E50") "  -0.27%  "

Set-TextValue $ws.Range("D51") "7.57"
Set-TextValue $ws.Range("E51") "  -2.06%  "
